$d = $word.ActiveDocument

# 1. Title text change: "Perceived Skill Gap Survey" -> "Survey Study on the Prestige of Alternative Education"
$d.Content.Find.Execute("Perceived Skill Gap Survey", $true, $false, $false, $false, $false, $true, 1, $false, "Survey Study on the Prestige of Alternative Education", 2) | Out-Null

# 2. Remove the bookmarks (the _Hlk... around the old title text, and the _GoBack one)
$d.Bookmarks.Item("_Hlk56671814").Delete()
$d.Bookmarks.Item("_GoBack").Delete()

# 3. Update the instructions paragraph: drop "Select the link below to complete the survey. " and
#    change "a code" -> "a completion code"
$d.Content.Find.Execute("We are conducting an academic survey about alternatives to the traditional college experience. Select the link below to complete the survey. At the end of the survey, you will receive a code to paste into the box below to receive credit for taking our survey. The survey will also ask for your Worker ID.", $true, $false, $false, $false, $false, $true, 1, $false, "We are conducting an academic survey about alternatives to the traditional college experience. At the end of the survey, you will receive a completion code to paste into the box below to receive credit for taking our survey. The survey will also ask for your Worker ID.", 2) | Out-Null

# 4. Replace the consent paragraph sentence with the new "To access the survey..." sentence
$d.Content.Find.Execute("As part of the survey, you will be asked to consent to participation in the study subject to an informed consent disclosure, which you may optionally preview now:", $true, $false, $false, $false, $false, $true, 1, $false, "To access the survey, click the button which says " + [char]8220 + "I agree" + [char]8221 + " at the bottom of the informed consent document below:", 2) | Out-Null

# 5. Replace the hyperlink paragraph entirely so it picks up clean (unformatted) pPr/rPr instead of
#    carrying over the old direct rFonts formatting, then re-insert a fresh hyperlink with the new URL.
$hyperlinkPara = $d.Paragraphs.Item(13)
$hyperlinkPara.Range.Delete()

$newPara = $d.Paragraphs.Item(13)
$newPara.Range.InsertBefore("https://osf.io/34rw8" + [char]13)

$refreshedPara = $d.Paragraphs.Item(13)
$urlRange = $d.Range($refreshedPara.Range.Start, $refreshedPara.Range.End - 1)
$d.Hyperlinks.Add($urlRange, "https://osf.io/34rw8") | Out-Null
